# Fruta / hortaliza, semanal
# A new weekly data row is inserted at row 29 (pushing the existing rows
# 29-64 down to 30-65), and the new row 29 is populated with this week's
# data for "Macroferia Regional de Talca" / Alcachofa / Madrigal / Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 29, shifting rows 29-64 down to 30-65.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new week's record.
$ws.Cells.Item(29, 1).Value = 5
$ws.Cells.Item(29, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(29, 3).Value = "Maule"
$ws.Cells.Item(29, 4).Value = 44484
$ws.Cells.Item(29, 5).Value = 7
$ws.Cells.Item(29, 6).Value = 100112013
$ws.Cells.Item(29, 7).Value = "Alcachofa"
$ws.Cells.Item(29, 8).Value = "Madrigal"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 300
$ws.Cells.Item(29, 11).Value = 10000
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 13).Value = 10000
$ws.Cells.Item(29, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(29, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(29, 16).Value = 250
$ws.Cells.Item(29, 17).Value = 40
$ws.Cells.Item(29, 18).Value = "Hortaliza"
